$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3000.0625
$ws.Range("I43").Value = 3083.3333
$ws.Range("J43").Value = 2750.25
$ws.Range("K43").Value = 3083.3333
$ws.Range("L43").Value = 2750.25
$ws.Range("M43").Value = -3014.3333
$ws.Range("N43").Value = -2888.25

$ws.Range("H51").Value = 4516.6665
$ws.Range("I51").Value = 3200
$ws.Range("J51").Value = 5833.3335
$ws.Range("K51").Value = 3200
$ws.Range("L51").Value = 5833.3335
$ws.Range("M51").Value = -2716
$ws.Range("N51").Value = -6801.3335

$ws.Range("H137").Value = 3347.6667
$ws.Range("I137").Value = 1423.2858
$ws.Range("J137").Value = 3564.9355
$ws.Range("K137").Value = 4269.857400000001
$ws.Range("L137").Value = 10694.8065
$ws.Range("M137").Value = -1719.857400000001
$ws.Range("N137").Value = -15794.8065

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29560.46
$ws.Range("I32").Value = 28324.285
$ws.Range("K32").Value = 28324.285
$ws.Range("M32").Value = -28037.285

$ws.Range("H74").Value = 1970.579
$ws.Range("I74").Value = 1169.2727
$ws.Range("J74").Value = 3072.375
$ws.Range("K74").Value = 1169.2727
$ws.Range("L74").Value = 3072.375
$ws.Range("M74").Value = -295.2727
$ws.Range("N74").Value = -4820.375

$ws.Range("H77").Value = 1970.579
$ws.Range("I77").Value = 1169.2727
$ws.Range("J77").Value = 3072.375
$ws.Range("K77").Value = 5846.363499999999
$ws.Range("L77").Value = 15361.875
$ws.Range("M77").Value = -1478.363499999999
$ws.Range("N77").Value = -24097.875

$ws.Range("H139").Value = 40598.375
$ws.Range("J139").Value = 40598.375
$ws.Range("L139").Value = 40598.375
$ws.Range("N139").Value = -50878.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 11108.889
$ws.Range("J81").Value = 11108.889
$ws.Range("L81").Value = 11108.889
$ws.Range("N81").Value = -13230.889

$ws.Range("H84").Value = 11108.889
$ws.Range("J84").Value = 11108.889
$ws.Range("L84").Value = 33326.667
$ws.Range("N84").Value = -43934.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 3500
$ws.Range("I54").Value = 3500
$ws.Range("K54").Value = 3500
$ws.Range("M54").Value = -2842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1513.8096
$ws.Range("J34").Value = 1822.9412
$ws.Range("L34").Value = 5468.8236
$ws.Range("N34").Value = -5636.8236

$ws.Range("H39").Value = 1293.6
$ws.Range("I39").Value = 1152
$ws.Range("J39").Value = 1315.3846
$ws.Range("K39").Value = 3456
$ws.Range("L39").Value = 3946.1538
$ws.Range("M39").Value = -3162
$ws.Range("N39").Value = -4534.1538

$ws.Range("H55").Value = 1220
$ws.Range("I55").Value = 700
$ws.Range("J55").Value = 1566.6666
$ws.Range("K55").Value = 2100
$ws.Range("L55").Value = 4699.9998
$ws.Range("M55").Value = -1923
$ws.Range("N55").Value = -5053.9998

$ws.Range("H63").Value = 4600
$ws.Range("J63").Value = 4600
$ws.Range("L63").Value = 13800
$ws.Range("N63").Value = -15298

$ws.Range("H66").Value = 4600
$ws.Range("J66").Value = 4600
$ws.Range("L66").Value = 41400
$ws.Range("N66").Value = -48888

$ws.Range("H69").Value = 38128570
$ws.Range("I69").Value = 1966.6666
$ws.Range("J69").Value = 42894396
$ws.Range("K69").Value = 5899.9998
$ws.Range("L69").Value = 128683188
$ws.Range("M69").Value = -5088.9998
$ws.Range("N69").Value = -128684810

$ws.Range("H72").Value = 38128570
$ws.Range("I72").Value = 1966.6666
$ws.Range("J72").Value = 42894396
$ws.Range("K72").Value = 17699.9994
$ws.Range("L72").Value = 386049564
$ws.Range("M72").Value = -13643.9994
$ws.Range("N72").Value = -386057676

$ws.Range("H104").Value = 3823.6667
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 3823.6667
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 11471.0001
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -16713.0001

$ws.Range("H109").Value = 12711.033
$ws.Range("I109").Value = 54563.4
$ws.Range("J109").Value = 4340.56
$ws.Range("K109").Value = 163690.2
$ws.Range("L109").Value = 13021.68
$ws.Range("M109").Value = -162650.2
$ws.Range("N109").Value = -15101.68

$ws.Range("H110").Value = 4000
$ws.Range("J110").Value = 6000
$ws.Range("L110").Value = 18000
$ws.Range("N110").Value = -26180

$ws.Range("H122").Value = 3025.6
$ws.Range("I122").Value = 704.0513
$ws.Range("J122").Value = 18115.666
$ws.Range("K122").Value = 6336.4617
$ws.Range("L122").Value = 163040.994
$ws.Range("M122").Value = -3886.4617
$ws.Range("N122").Value = -167940.994

$ws.Range("H132").Value = 1524.625
$ws.Range("I132").Value = 1207.8334
$ws.Range("J132").Value = 2475
$ws.Range("K132").Value = 10870.5006
$ws.Range("L132").Value = 22275
$ws.Range("M132").Value = -8340.500599999999
$ws.Range("N132").Value = -27335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 794.2857
$ws.Range("I12").Value = 794.2857
$ws.Range("K12").Value = 794.2857
$ws.Range("M12").Value = -654.2857

$ws.Range("H138").Value = 42272.727
$ws.Range("J138").Value = 42272.727
$ws.Range("L138").Value = 42272.727
$ws.Range("N138").Value = -52552.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2411.647
$ws.Range("I40").Value = 2321.3572
$ws.Range("J40").Value = 2833
$ws.Range("K40").Value = 2321.3572
$ws.Range("L40").Value = 2833
$ws.Range("M40").Value = -2185.3572
$ws.Range("N40").Value = -3105

$ws.Range("H54").Value = 14997.5
$ws.Range("J54").Value = 14997.5
$ws.Range("L54").Value = 14997.5
$ws.Range("N54").Value = -16285.5

$ws.Range("H132").Value = 4077.8538
$ws.Range("I132").Value = 3938.652
$ws.Range("J132").Value = 4255.722
$ws.Range("K132").Value = 11815.956
$ws.Range("L132").Value = 12767.166
$ws.Range("M132").Value = -9285.956
$ws.Range("N132").Value = -17827.166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 76670.664
$ws.Range("J19").Value = 76670.664
$ws.Range("L19").Value = 76670.664
$ws.Range("N19").Value = -77018.664

$ws.Range("H46").Value = 62399.08
$ws.Range("J46").Value = 62399.08
$ws.Range("L46").Value = 62399.08
$ws.Range("N46").Value = -62861.08

$ws.Range("H134").Value = 62399.08
$ws.Range("J134").Value = 62399.08
$ws.Range("L134").Value = 187197.24
$ws.Range("N134").Value = -192267.24
